# Apply cryptocurrency price/volume updates (GitHub Actions cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.0000", "27.489.38") that must
# stay literal text rather than being auto-converted to numbers by Excel, so
# force the Text number format on the whole data range up front.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.489.38"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "1.816.39"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "344.27"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "0.3830"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").Value = "0.3551"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("D9").Value = "49.02"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "1.236"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("E11").Value = "  +3.85%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "22.42"
$ws.Range("E13").Value = "  +9.29%  "
$ws.Range("D14").Value = "6.610"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "1.815.07"
$ws.Range("E15").Value = "  +3.98%  "
$ws.Range("D16").Value = "7.228"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "0.00001129"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("D18").Value = "0.06734"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "86.95"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "17.66"
$ws.Range("E21").Value = "  +5.02%  "
$ws.Range("D22").Value = "6.557"
$ws.Range("E22").Value = "  +5.73%  "
$ws.Range("D23").Value = "13.15"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "27.489.14"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("D25").Value = "2.475"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "2.702"
$ws.Range("E26").Value = "  +7.21%  "
$ws.Range("D27").Value = "22.20"
$ws.Range("E27").Value = "  +12.74%  "
$ws.Range("D28").Value = "1.468"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "154.03"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "2.019.80"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").Value = "135.90"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").Value = "6.398"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("D33").Value = "4.060"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").Value = "13.93"
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("D35").Value = "0.08821"
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").Value = "1.698"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "5.647"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("D38").Value = "0.7021"
$ws.Range("E38").Value = "  +11.71%  "
$ws.Range("D39").Value = "9.069"
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06516"
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.02411"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "0.2257"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("D43").Value = "1.301"
$ws.Range("E43").Value = "  +5.30%  "
$ws.Range("D44").Value = "14.95"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("E45").Value = "  +8.59%  "
$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "3.972"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "2.193"
$ws.Range("E48").Value = "  +5.74%  "
$ws.Range("D49").Value = "133.28"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Value = "0.07327"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "80.94"
$ws.Range("E51").Value = "  +3.74%  "
